# Update results of using trust-region method to solve continuous
# relaxation with TV regularizer.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new rows at position 4 (GRAPE+TR / ADMM+TR continuous
#    relaxation results), pushing the old rows 4-20 down to 6-22.
$ws.Rows("4:5").Insert()

# 2. Fill in new row 4: GRAPE+TR (continuous)
$ws.Range("A4").Value = ""
$ws.Range("B4").Value = "Energy2_evotime2.0_n_ts40_ptypeCONSTANT_offset0.5_sigma0.25_eta0.001_threshold30_iter100_typetvc"
$ws.Range("C4").Value = "GRAPE+TR (continuous)"
$ws.Range("D4").Value = -0.999
$ws.Range("E4").Value = 0.567
$ws.Range("F4").Formula = "=D4+E4*0.01"
$ws.Range("G4").Value = 1.77
$ws.Range("H4").Value = "solve continuous relaxation`nwith TV norm starting from GRAPE"

# 3. Fill in new row 5: ADMM+TR (continuous)
$ws.Range("B5").Value = "EnergyADMM2_evotime2.0_n_ts40_ptypeWARM_offset0.5_penalty0.01_ADMM_10.0_iter100_sigma0.25_eta0.001_threshold30_iter100_typetvc"
$ws.Range("C5").Value = "ADMM+TR (continuous)"
$ws.Range("D5").Value = -0.999
$ws.Range("E5").Value = 0.523
$ws.Range("F5").Formula = "=D5+E5*0.01"
$ws.Range("G5").Value = 0.09
$ws.Range("H5").Value = "solve continuous relaxation`nwith TV norm starting from ADMM"

# 4. Relabel the legacy "*+TR" rows to "*+TR (binary)" -- they now live
#    at rows 14, 17, 20 after the insert.
$ws.Range("C14").Value = "GRAPE+TR (binary)"
$ws.Range("C17").Value = "ADMM+TR (binary)"
$ws.Range("C20").Value = "Switching+TR (binary)"

# 5. Header G1: was "time (s)" already but shared-string table shrank so
#    just make sure it is still correct.
$ws.Range("G1").Value = "time (s)"

# 6. Append the SNOPT / IPOPT continuous-result rows at the bottom.
$ws.Range("C25").Value = "SNOPT"
$ws.Range("D25").Value = -0.916
$ws.Range("H25").Value = "continuous results"

$ws.Range("C26").Value = "IPOPT"
$ws.Range("D26").Value = -0.916
$ws.Range("G26").Value = 0.02
$ws.Range("H26").Value = "continuous results"

# 7. Column width / view tweaks captured in the diff.
$ws.Columns("B").ColumnWidth = 120.5
$ws.Columns("C").ColumnWidth = 23
$ws.Columns("H").ColumnWidth = 31.6640625

$ws.Range("A14").Select()
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Range("C31").Select()
